$wb = $excel.ActiveWorkbook

# OFF sheet - Road ("R") row - Week 17 data
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 234
$wsOff.Range("C3").Value = 176
$wsOff.Range("D3").Value = 61
$wsOff.Range("E3").Value = 36

# DEF sheet - Road ("R") row - Week 17 data
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 198
$wsDef.Range("C3").Value = 133
$wsDef.Range("D3").Value = 35
$wsDef.Range("E3").Value = 17
$wsDef.Range("F3").Value = 6
